$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (column D) values; force text format so strings like
# "36.10" / "9.80" / "1.00" keep their literal trailing zeros instead of
# being reinterpreted as numbers. ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.002.98"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.841.02"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "636.68"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.81"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.843.15"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.66"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000253"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.10"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.482.18"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.897.40"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.887.02"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.17"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.17"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "470.51"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.80"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.711"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000152"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.90"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.04"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.12"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.987.59"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.35"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.40"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.782.58"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.45"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.94"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.985"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "157.55"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.95"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "47.41"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.42"

# --- Update Volume(1h) (column E) values ---
$ws.Range("E2").Value = "  +2.85%  "
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("E5").Value = "  +6.00%  "
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("E7").Value = "  +1.40%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("E10").Value = "  +2.40%  "
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("E12").Value = "  +2.90%  "
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("E16").Value = "  +2.77%  "
$ws.Range("E17").Value = "  +2.71%  "
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("E21").Value = "  +2.02%  "
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("E23").Value = "  +1.67%  "
$ws.Range("E24").Value = "  +2.16%  "
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("E26").Value = "  +3.74%  "
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("E31").Value = "  +2.96%  "
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("E35").Value = "  +1.09%  "
$ws.Range("E38").Value = "  +3.27%  "
$ws.Range("E39").Value = "  +8.61%  "
$ws.Range("E40").Value = "  +5.26%  "
$ws.Range("E41").Value = "  +2.72%  "
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("E45").Value = "  +3.93%  "
$ws.Range("E46").Value = "  +1.70%  "
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("E49").Value = "  +3.62%  "
$ws.Range("E50").Value = "  +5.11%  "
$ws.Range("E51").Value = "  +1.50%  "

# --- Row 36/37 reorder: Binance-PegBSC-USD <-> Aptos swap places ---
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.10"
$ws.Range("E36").Value = "  +1.16%  "
$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.02%  "
